$wb = $excel.ActiveWorkbook

# --- Folder Inventory sheet ---
$ws = $wb.Worksheets.Item("Folder Inventory")

# Insert a brand new row at position 2 (pushes existing rows down by one)
$ws.Rows.Item(2).Insert()
# The inserted row inherits the header's bold/centered style; strip it back to plain
$ws.Range("A2:E2").ClearFormats()

$ws.Range("A2").Value = "Power_Platform_Workshop:Administration_and_Governance"
$ws.Range("B2").Value = "Power_Platform_Workshop:Administration_and_Governance"
$ws.Range("C2").Value = "2025-06-16 10:53:14 +0530"
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = "Root"

# The old stale duplicate entry for this same folder (now shifted to row 50) is removed
$ws.Rows.Item(50).Delete()

# --- Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "2025-06-16 05:23:34 UTC"
$meta.Range("B5").NumberFormat = "@"
$meta.Range("B5").Value = "8"

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B5").Value = "2025-06-16 10:53:14 +0530"
